$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Fgf1","Cd44","ECs",2,0.6666666666666666,0.265349,0.7960469999999999,0.1498685997319469,0.1498685997319469,3,1,19.21315233333334,57.63945700000001,0.04451179209991234,0.04451179209991233,5.098190758497668,45.883716826479,0.006670919953573396,0.006670919953573396),
    @("ECs","Fgf1","Cd44","FAPs",2,0.6666666666666666,0.265349,0.7960469999999999,0.1498685997319469,0.1498685997319469,3,1,92.44713066666667,277.341392,0.2141755495962477,0.2141755495962477,24.53075367526933,220.776783077424,0.03209818971480979,0.03209818971480978),
    @("ECs","Fgf1","Cd44","Inflammatory-Mac",2,0.6666666666666666,0.265349,0.7960469999999999,0.1498685997319469,0.1498685997319469,3,1,166.8580016666666,500.5740049999999,0.3865658561145097,0.3865658561145097,44.27560388424833,398.4804349582349,0.05793408356006282,0.05793408356006282),
    @("ECs","Fgf1","Cd44","MuSCs",2,0.6666666666666666,0.265349,0.7960469999999999,0.1498685997319469,0.1498685997319469,3,1,41.09915599999999,123.297468,0.09521587377309249,0.09521587377309249,10.905619945444,98.15057950899597,0.01426986967462717,0.01426986967462717),
    @("ECs","Fgf1","Cd44","Resolving-Mac",2,0.6666666666666666,0.265349,0.7960469999999999,0.1498685997319469,0.1498685997319469,3,1,112.0244103333333,336.073231,0.2595309284162377,0.2595309284162377,29.72556525753967,267.530087317857,0.03889553682887368,0.03889553682887368),
    @("FAPs","Fgf1","Cd44","ECs",3,1,1.003400666666667,3.010202,0.5667187473230926,0.5667187473230925,3,1,19.21315233333334,57.63945700000001,0.04451179209991234,0.04451179209991233,19.27848986003489,173.506408740314,0.02522566705996825,0.02522566705996824),
    @("FAPs","Fgf1","Cd44","FAPs",3,1,1.003400666666667,3.010202,0.5667187473230926,0.5667187473230925,3,1,92.44713066666667,277.341392,0.2141755495962477,0.2141755495962477,92.76151254235377,834.8536128811838,0.1213772991744204,0.1213772991744204),
    @("FAPs","Fgf1","Cd44","Inflammatory-Mac",3,1,1.003400666666667,3.010202,0.5667187473230926,0.5667187473230925,3,1,166.8580016666666,500.5740049999999,0.3865658561145097,0.3865658561145097,167.4254301110011,1506.82887099901,0.2190741177350938,0.2190741177350937),
    @("FAPs","Fgf1","Cd44","MuSCs",3,1,1.003400666666667,3.010202,0.5667187473230926,0.5667187473230925,3,1,41.09915599999999,123.297468,0.09521587377309249,0.09521587377309249,41.23892052983732,371.1502847685359,0.05396062070996067,0.05396062070996067),
    @("FAPs","Fgf1","Cd44","Resolving-Mac",3,1,1.003400666666667,3.010202,0.5667187473230926,0.5667187473230925,3,1,112.0244103333333,336.073231,0.2595309284162377,0.2595309284162377,112.4053680114069,1011.648312102662,0.1470810426436494,0.1470810426436494),
    @("MuSCs","Fgf1","Cd44","ECs",3,1,0.4790693333333333,1.437208,0.2705774288246194,0.2705774288246194,3,1,19.21315233333334,57.63945700000001,0.04451179209991234,0.04451179209991233,9.204432079561778,82.839888716056,0.01204388625877029,0.01204388625877029),
    @("MuSCs","Fgf1","Cd44","FAPs",3,1,0.4790693333333333,1.437208,0.2705774288246194,0.2705774288246194,3,1,92.44713066666667,277.341392,0.2141755495962477,0.2141755495962477,44.28858525705955,398.597267313536,0.05795106952685247,0.05795106952685247),
    @("MuSCs","Fgf1","Cd44","Inflammatory-Mac",3,1,0.4790693333333333,1.437208,0.2705774288246194,0.2705774288246194,3,1,166.8580016666666,500.5740049999999,0.3865658561145097,0.3865658561145097,79.93655161978221,719.4289645780399,0.1045959954188518,0.1045959954188518),
    @("MuSCs","Fgf1","Cd44","MuSCs",3,1,0.4790693333333333,1.437208,0.2705774288246194,0.2705774288246194,3,1,41.09915599999999,123.297468,0.09521587377309249,0.09521587377309249,19.68934526548266,177.204107389344,0.02576326630881288,0.02576326630881288),
    @("MuSCs","Fgf1","Cd44","Resolving-Mac",3,1,0.4790693333333333,1.437208,0.2705774288246194,0.2705774288246194,3,1,112.0244103333333,336.073231,0.2595309284162377,0.2595309284162377,53.66745957544978,483.007136179048,0.07022321131133195,0.07022321131133195),
    @("Resolving-Mac","Fgf1","Cd44","ECs",1,0.3333333333333333,0.02272533333333333,0.068176,0.01283522412034115,0.01283522412034115,3,1,19.21315233333334,57.63945700000001,0.04451179209991234,0.04451179209991233,0.4366252911591112,3.929627620432,0.0005713188276004053,0.0005713188276004052),
    @("Resolving-Mac","Fgf1","Cd44","FAPs",1,0.3333333333333333,0.02272533333333333,0.068176,0.01283522412034115,0.01283522412034115,3,1,92.44713066666667,277.341392,0.2141755495962477,0.2141755495962477,2.100891860110222,18.908026740992,0.00274899118016508,0.00274899118016508),
    @("Resolving-Mac","Fgf1","Cd44","Inflammatory-Mac",1,0.3333333333333333,0.02272533333333333,0.068176,0.01283522412034115,0.01283522412034115,3,1,166.8580016666666,500.5740049999999,0.3865658561145097,0.3865658561145097,3.791903707208888,34.12713336488,0.00496165940050128,0.00496165940050128),
    @("Resolving-Mac","Fgf1","Cd44","MuSCs",1,0.3333333333333333,0.02272533333333333,0.068176,0.01283522412034115,0.01283522412034115,3,1,41.09915599999999,123.297468,0.09521587377309249,0.09521587377309249,0.9339920198186665,8.405928178367999,0.001222117079691754,0.001222117079691754),
    @("Resolving-Mac","Fgf1","Cd44","Resolving-Mac",1,0.3333333333333333,0.02272533333333333,0.068176,0.01283522412034115,0.01283522412034115,3,1,112.0244103333333,336.073231,0.2595309284162377,0.2595309284162377,2.545792066295111,22.912128596656,0.003331137632382626,0.003331137632382626)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowValues[$j]
    }
}

Write-Host "Updated $($data.Count) rows"

